# Insert a new weekly price record for "Apio" (Macroferia Regional de Talca)
# just above the former row 111. Excel shifts the existing rows 111-219
# down to 112-220 and the new row inherits the same Mercado/Categoría
# metadata as the row that used to occupy position 111, while getting a
# fresh date (2022-08-17) and volume (700).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: push rows 111..219 down to 112..220.
$ws.Range("A111:R111").Insert()

# Seed the new row with the same values/format the (now shifted) row 112
# holds -- this is what used to be row 111 before the insert -- so every
# column except the ones that actually change keeps its data intact.
$ws.Range("A112:R112").Copy($ws.Range("A111:R111"))

# Apply the two values that differ for this new record.
$ws.Range("D111").Value = 44790
$ws.Range("J111").Value = 700
